# exercise_02_Sol.xlsx — update the "start of year 20" annuity PV note to
# "start of year 21", and add a cross-check formula for that PV using the
# direct annuity-difference method.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Solutions")

# 1) A183: "The PV of this annuity at the start of year 20 is"
#          -> "The PV of this annuity at the start of year 21 is"
$ws.Range("A183").Value = "The PV of this annuity at the start of year 21 is"

# 2) D186: new cross-check formula computing the same PV a different way
$ws.Range("D186").Formula = "=10000/0.05*(1 - 1/(1.05)^45)-10000/0.05*(1 - 1/(1.05)^21)"

# 3) Row 170 had an explicit wrapped-text row height; auto-fit it back to
#    the sheet's default height.
$ws.Rows.Item(170).AutoFit()

# 4) Leave the view scrolled down to the annuity section with B186 selected,
#    matching where the edit was made.
$ws.Activate() | Out-Null
$ws.Range("B186").Select() | Out-Null
